# Dillards ROI.xlsx - "Add files via upload" edit
#
# Author lowered the "% Stores Successful" assumption (B28) from 26.46% to
# 15.7%. B30 (Proj. Successful Stores = B16*B28) and B32 (Annual Profit /
# store = B30*B25) are formulas and recalculate automatically.
#
# The file was also re-saved with the window scrolled down (top-left cell
# A4) with gridlines turned off, and the last selection left on D25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core data edit -------------------------------------------------------
# % Stores Successful: 26.46% -> 15.7%
$ws.Range("B28").Value = 0.157

# --- View / window cosmetics ---------------------------------------------
# Hide gridlines on the sheet.
$excel.ActiveWindow.DisplayGridlines = $false

# Scroll so row 4 is the top visible row (sheetView topLeftCell="A4").
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# Leave the final selection on D25.
$ws.Range("D25").Select()
